$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 884257.4399999999
$ws.Range("I17").Value = 2437.5386
$ws.Range("J17").Value = 2317214.8
$ws.Range("K17").Value = 7312.6158
$ws.Range("L17").Value = 6951644.399999999
$ws.Range("M17").Value = -7144.6158
$ws.Range("N17").Value = -6951980.399999999
# Row 32
$ws.Range("H32").Value = 9896
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 9896
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 9896
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10548
# Row 33
$ws.Range("H33").Value = 4348182
$ws.Range("I33").Value = 5000358
$ws.Range("J33").Value = 340.66666
$ws.Range("K33").Value = 5000358
$ws.Range("L33").Value = 340.66666
$ws.Range("M33").Value = -5000129
$ws.Range("N33").Value = -798.66666
# Row 41
$ws.Range("H41").Value = 541.8095
$ws.Range("I41").Value = 553.75
$ws.Range("J41").Value = 525.8889
$ws.Range("K41").Value = 553.75
$ws.Range("L41").Value = 525.8889
$ws.Range("M41").Value = -113.75
$ws.Range("N41").Value = -1405.8889
# Row 64
$ws.Range("H64").Value = 8808.5
$ws.Range("J64").Value = 10135.625
$ws.Range("L64").Value = 10135.625
$ws.Range("N64").Value = -10631.625
# Row 67
$ws.Range("H67").Value = 8808.5
$ws.Range("J67").Value = 10135.625
$ws.Range("L67").Value = 10135.625
$ws.Range("N67").Value = -11851.625
# Row 98
$ws.Range("H98").Value = 1184.875
$ws.Range("I98").Value = 1279.8334
$ws.Range("K98").Value = 1279.8334
$ws.Range("M98").Value = 218.1666
# Row 106
$ws.Range("H106").Value = 1927.8667
$ws.Range("I106").Value = 1355.3636
$ws.Range("K106").Value = 1355.3636
$ws.Range("M106").Value = -724.3635999999999
# Row 111
$ws.Range("H111").Value = 2832.0715
$ws.Range("I111").Value = 527.375
$ws.Range("J111").Value = 5905
$ws.Range("K111").Value = 1582.125
$ws.Range("L111").Value = 17715
$ws.Range("M111").Value = 1484.875
$ws.Range("N111").Value = -23849
# Row 113
$ws.Range("H113").Value = 14663.77
$ws.Range("I113").Value = 12539.667
$ws.Range("K113").Value = 12539.667
$ws.Range("M113").Value = -9285.666999999999
# Row 114
$ws.Range("H114").Value = 45000
$ws.Range("J114").Value = 45000
$ws.Range("L114").Value = 45000
$ws.Range("N114").Value = -53678
# Row 122
$ws.Range("H122").Value = 1184.875
$ws.Range("I122").Value = 1279.8334
$ws.Range("K122").Value = 3839.5002
$ws.Range("M122").Value = -1389.5002
# Row 131
$ws.Range("H131").Value = 13809.533
$ws.Range("I131").Value = 3782.4443
$ws.Range("J131").Value = 28850.166
$ws.Range("K131").Value = 11347.3329
$ws.Range("L131").Value = 86550.49800000001
$ws.Range("M131").Value = -6307.332900000001
$ws.Range("N131").Value = -96630.49800000001
# Row 137
$ws.Range("H137").Value = 2589.4
$ws.Range("I137").Value = 3079.2
$ws.Range("J137").Value = 2099.6
$ws.Range("K137").Value = 9237.599999999999
$ws.Range("L137").Value = 6298.799999999999
$ws.Range("M137").Value = -6687.599999999999
$ws.Range("N137").Value = -11398.8
# Row 141
$ws.Range("H141").Value = 2096.0625
$ws.Range("I141").Value = 2139.1333
$ws.Range("K141").Value = 6417.3999
$ws.Range("M141").Value = -1237.3999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1947.8667
$ws.Range("I2").Value = 1939.9231
$ws.Range("J2").Value = 1999.5
$ws.Range("K2").Value = 1939.9231
$ws.Range("L2").Value = 1999.5
$ws.Range("M2").Value = -1826.9231
$ws.Range("N2").Value = -2225.5
# Row 22
$ws.Range("H22").Value = 5199.75
$ws.Range("J22").Value = 9899
$ws.Range("L22").Value = 9899
$ws.Range("N22").Value = -10497
# Row 32
$ws.Range("H32").Value = 3042.8281
$ws.Range("I32").Value = 2948.2698
$ws.Range("K32").Value = 2948.2698
$ws.Range("M32").Value = -2661.2698
# Row 63
$ws.Range("H63").Value = 2795.1
$ws.Range("J63").Value = 2299
$ws.Range("L63").Value = 2299
$ws.Range("N63").Value = -3671
# Row 66
$ws.Range("H66").Value = 2795.1
$ws.Range("J66").Value = 2299
$ws.Range("L66").Value = 11495
$ws.Range("N66").Value = -18359
# Row 116
$ws.Range("H116").Value = 1947.8667
$ws.Range("I116").Value = 1939.9231
$ws.Range("J116").Value = 1999.5
$ws.Range("K116").Value = 1939.9231
$ws.Range("L116").Value = 1999.5
$ws.Range("M116").Value = 354.0769
$ws.Range("N116").Value = -6587.5
# Row 132
$ws.Range("H132").Value = 5508.5938
$ws.Range("I132").Value = 5171.846
$ws.Range("J132").Value = 6967.8335
$ws.Range("K132").Value = 15515.538
$ws.Range("L132").Value = 20903.5005
$ws.Range("M132").Value = -12985.538
$ws.Range("N132").Value = -25963.5005

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1947.8667
$ws.Range("I3").Value = 1939.9231
$ws.Range("J3").Value = 1999.5
$ws.Range("K3").Value = 1939.9231
$ws.Range("L3").Value = 1999.5
$ws.Range("M3").Value = -1825.9231
$ws.Range("N3").Value = -2227.5
# Row 20
$ws.Range("H20").Value = 3231.4102
$ws.Range("I20").Value = 2800.261
$ws.Range("J20").Value = 3851.1875
$ws.Range("K20").Value = 2800.261
$ws.Range("L20").Value = 3851.1875
$ws.Range("M20").Value = -2553.261
$ws.Range("N20").Value = -4345.1875
# Row 99
$ws.Range("H99").Value = 3757.7368
$ws.Range("I99").Value = 1846.8823
$ws.Range("K99").Value = 1846.8823
$ws.Range("M99").Value = -348.8823
# Row 134
$ws.Range("H134").Value = 1845.4576
$ws.Range("I134").Value = 1837.1786
$ws.Range("K134").Value = 5511.5358
$ws.Range("M134").Value = -2976.5358

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 1304.909
$ws.Range("I7").Value = 1600.25
$ws.Range("K7").Value = 1600.25
$ws.Range("M7").Value = -1487.25
# Row 70
$ws.Range("H70").Value = 55000
$ws.Range("J70").Value = 55000
$ws.Range("L70").Value = 55000
$ws.Range("N70").Value = -55630
# Row 73
$ws.Range("H73").Value = 55000
$ws.Range("J73").Value = 55000
$ws.Range("L73").Value = 55000
$ws.Range("N73").Value = -57184

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 64375750
$ws.Range("J4").Value = 171
$ws.Range("L4").Value = 513
$ws.Range("N4").Value = -737
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 113
$ws.Range("H113").Value = 694.6111
$ws.Range("I113").Value = 775.1429000000001
$ws.Range("J113").Value = 643.36365
$ws.Range("K113").Value = 2325.4287
$ws.Range("L113").Value = 1930.09095
$ws.Range("M113").Value = -155.4287000000004
$ws.Range("N113").Value = -6270.09095
# Row 117
$ws.Range("H117").Value = 854.75
$ws.Range("J117").Value = 884.7143
$ws.Range("L117").Value = 2654.1429
$ws.Range("N117").Value = -9538.142899999999
# Row 120
$ws.Range("H120").Value = 9447.5
$ws.Range("I120").Value = 9447.5
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 28342.5
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -23504.5
$ws.Range("N120").ClearContents()
# Row 121
$ws.Range("H121").Value = 1483.6364
$ws.Range("I121").Value = 2385.6
$ws.Range("J121").Value = 732
$ws.Range("K121").Value = 7156.799999999999
$ws.Range("L121").Value = 2196
$ws.Range("M121").Value = -5846.799999999999
$ws.Range("N121").Value = -4816
# Row 129
$ws.Range("H129").Value = 1476.5555
$ws.Range("I129").Value = 1041.2858
$ws.Range("K129").Value = 3123.8574
$ws.Range("M129").Value = 1876.1426

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 139.60869
$ws.Range("I2").Value = 85.15000000000001
$ws.Range("J2").Value = 502.66666
$ws.Range("K2").Value = 85.15000000000001
$ws.Range("L2").Value = 502.66666
$ws.Range("M2").Value = 27.84999999999999
$ws.Range("N2").Value = -728.66666
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 118
$ws.Range("H118").Value = 23875
$ws.Range("J118").Value = 23875
$ws.Range("L118").Value = 23875
$ws.Range("N118").Value = -27189
# Row 132
$ws.Range("H132").Value = 6274.3335
$ws.Range("I132").Value = 6476.8237
$ws.Range("J132").Value = 4897.4
$ws.Range("K132").Value = 19430.4711
$ws.Range("L132").Value = 14692.2
$ws.Range("M132").Value = -16900.4711
$ws.Range("N132").Value = -19752.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6380.7144
$ws.Range("I7").Value = 6156.625
$ws.Range("K7").Value = 6156.625
$ws.Range("M7").Value = -6044.625
# Row 40
$ws.Range("H40").Value = 6403.722
$ws.Range("I40").Value = 5734.5835
$ws.Range("K40").Value = 5734.5835
$ws.Range("M40").Value = -5598.5835
# Row 46
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376
# Row 122
$ws.Range("H122").Value = 481814.75
$ws.Range("I122").Value = 718886.7
$ws.Range("J122").Value = 7670.857
$ws.Range("K122").Value = 2156660.1
$ws.Range("L122").Value = 23012.571
$ws.Range("M122").Value = -2154210.1
$ws.Range("N122").Value = -27912.571
# Row 126
$ws.Range("H126").Value = 6380.7144
$ws.Range("I126").Value = 6156.625
$ws.Range("K126").Value = 18469.875
$ws.Range("M126").Value = -15999.875
